# Apply the "Added 'we're hiring' slide to end of presentation" edit.
#
# 1) Merge the split "Thanks Curtissimo" / "!" runs on the title slide into
#    a single run reading "Thanks Curtissimo!".
# 2) Append two new "Title and Content" slides at the end of the deck:
#      - "Questions And Answers"
#      - "We're Hiring"

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Title slide: merge "Thanks Curtissimo" + "!" runs into one run.
# ---------------------------------------------------------------------
$titleSlide = $p.Slides.Item(1)
$titleShape = $titleSlide.Shapes.Item(1)
$titleRange = $titleShape.TextFrame.TextRange

$fullTitleText = $titleRange.Text
$mergedText = "Thanks Curtissimo!"
$startPos = $fullTitleText.IndexOf("Thanks Curtissimo") + 1

$firstRun = $titleRange.Characters($startPos, 17)
$firstRun.Text = $mergedText

$trailingBang = $titleRange.Characters($titleRange.Length, 1)
$trailingBang.Text = ""

# ---------------------------------------------------------------------
# 2) Append the two new slides using the same "Title and Content" layout
#    as the rest of the deck's content slides.
# ---------------------------------------------------------------------
$contentLayout = $p.Slides.Item($p.Slides.Count).CustomLayout

# --- Slide: "Questions And Answers" --------------------------------
$qaSlide = $p.Slides.AddSlide($p.Slides.Count + 1, $contentLayout)

$qaSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Questions And Answers"

$qaBody = $qaSlide.Shapes.Item(2).TextFrame.TextRange
$qaBody.Text = "You line " + [char]0x2018 + "em up." + "`r" + "I" + [char]0x2019 + "ll knock " + [char]0x2018 + "em down." + "`r" + "If I say something dumb and you know better, SPEAK UP!"

# split out "em" in the first line so it mirrors the authored run layout
$p1 = $qaBody.Paragraphs(1, 1)
$emIdx = $p1.Text.IndexOf("em") + 1
$p1.Characters($emIdx, 2).Text = "em"

# split out "em" in the second line as well
$p2 = $qaBody.Paragraphs(2, 1)
$emIdx2 = $p2.Text.IndexOf("em") + 1
$p2.Characters($emIdx2, 2).Text = "em"

# bold + red for the third line
$p3 = $qaBody.Paragraphs(3, 1)
$p3.Font.Bold = $true
$p3.Font.Color.RGB = 255

# --- Slide: "We're Hiring" ------------------------------------------
$hiringSlide = $p.Slides.AddSlide($p.Slides.Count + 1, $contentLayout)

$hiringSlide.Shapes.Item(1).TextFrame.TextRange.Text = "We" + [char]0x2019 + "re Hiring"

$hiringBody = $hiringSlide.Shapes.Item(2).TextFrame.TextRange
$hiringBody.Text = "Senior .Net Developer`rBusiness Development Manager`rProgram Manager`r`rEmail Valerie Carmona:  Valerie.Carmona@improvingenterprises.com"

# split out ".Net" in the first line so it mirrors the authored run layout
$h1 = $hiringBody.Paragraphs(1, 1)
$netIdx = $h1.Text.IndexOf(".Net") + 1
$h1.Characters($netIdx, 4).Text = ".Net"

# last paragraph: no bullet / no indent, smaller font, with a mailto hyperlink
# on the email address
$h5 = $hiringBody.Paragraphs(5, 1)
$h5.ParagraphFormat.Bullet.Visible = $false
$h5.IndentLevel = 1
$h5.Font.Size = 20

$emailStart = $h5.Text.IndexOf("Valerie.Carmona@improvingenterprises.com") + 1
$emailLen = [string]"Valerie.Carmona@improvingenterprises.com".Length
$emailRun = $h5.Characters($emailStart, $emailLen)
$emailRun.ActionSettings.Item(1).Hyperlink.Address = "mailto:Valerie.Carmona@improvingenterprises.com"
